$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.902.47"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "3.163.39"
$ws.Range("E3").Value = "  -7.69%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.42%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.39%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.162.13"
$ws.Range("E9").Value = "  -7.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.13%  "
$ws.Range("E11").Value = "  -6.14%  "
$ws.Range("E12").Value = "  -5.00%  "
$ws.Range("D13").Value = "3.709.69"
$ws.Range("E13").Value = "  -7.85%  "
$ws.Range("E14").Value = "  +1.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.68%  "
$ws.Range("D16").Value = "64.878.84"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("E17").Value = "  -6.54%  "
$ws.Range("D18").Value = "3.165.13"
$ws.Range("E18").Value = "  -7.36%  "
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.52%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.497"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000116"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.176"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -4.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.72%  "
$ws.Range("E33").Value = "  -8.47%  "
$ws.Range("E34").Value = "  -5.04%  "
$ws.Range("E35").Value = "  -6.19%  "
$ws.Range("E36").Value = "  -6.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "155.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.89%  "
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.79%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.661.48"
$ws.Range("E41").Value = "  -1.73%  "
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.47"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.14%  "
$ws.Range("E43").Value = "  -5.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.04"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0657"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "325.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0273"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("E51").Value = "  -0.08%  "
